$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (knn)
$ws.Range("B4").Value = 0.5622630665831265
$ws.Range("C4").Value = 0.5599999999999999
$ws.Range("D4").Value = 0.5718680981367751
$ws.Range("E4").Value = 0.5625
$ws.Range("F4").Value = 0.4973900991999861
$ws.Range("G4").Value = 0.488
$ws.Range("H4").Value = 0.5115670629662663
$ws.Range("I4").Value = 0.5119999999999999
$ws.Range("J4").Value = 0.6478652111836987
$ws.Range("K4").Value = 0.671
$ws.Range("L4").Value = 0.6320553432887318
$ws.Range("M4").Value = 0.641

# Row 5 (nb)
$ws.Range("B5").Value = 0.6019897392385534
$ws.Range("C5").Value = 0.635
$ws.Range("D5").Value = 0.5994111323969925
$ws.Range("E5").Value = 0.5985
$ws.Range("F5").Value = 0.6586340925147478
$ws.Range("G5").Value = 0.958
$ws.Range("H5").Value = 0.5019838417455367
$ws.Range("I5").Value = 0.5044999999999999
$ws.Range("J5").Value = 0.6588398290076032
$ws.Range("K5").Value = 0.806
$ws.Range("L5").Value = 0.5635210902082279
$ws.Range("M5").Value = 0.5854999999999999

# Row 6 (rf)
$ws.Range("B6").Value = 0.5545416489452318
$ws.Range("C6").Value = 0.5539999999999999
$ws.Range("D6").Value = 0.5585990932195147
$ws.Range("E6").Value = 0.5555000000000001
$ws.Range("F6").Value = 0.4973844107952275
$ws.Range("G6").Value = 0.4940000000000001
$ws.Range("H6").Value = 0.5028591809916632
$ws.Range("I6").Value = 0.5035000000000001
$ws.Range("J6").Value = 0.6463574602214355
$ws.Range("K6").Value = 0.642
$ws.Range("L6").Value = 0.6627215819586342
$ws.Range("M6").Value = 0.6575
